$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell value is written with a leading apostrophe to force Excel to
# store it as literal text (otherwise numeric-looking strings such as
# "96.89" get silently converted to a number and lose their exact
# formatting / trailing zeros / thousands separators). Resetting the
# range style back to "Normal" afterwards clears the "quote prefix" cell
# style that the apostrophe trick leaves behind, so the saved cell keeps
# the same (default) style as before the edit.

$ws.Range("D2").Value = "'46.300.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.00%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.463.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +8.38%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'297.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.60%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'96.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.56%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +1.10%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.04%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.515"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.68%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'35.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.24%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0789"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.42%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'7.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.33%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +1.95%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.839.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +8.40%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.466.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +8.30%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +8.09%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  +4.34%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'46.333.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.76%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'12.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.63%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0951"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.70%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +8.30%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'67.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +3.06%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'246.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.11%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'1.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +6.93%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.06%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'40.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.47%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.58%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +3.57%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'3.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +16.57%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'21.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +7.66%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "'WEMIXToken"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'2.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.26%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = "'Filecoin"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'5.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +5.43%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'148.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.55%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +23.68%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.0778"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.82%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +2.44%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.78%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'15.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.16%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +4.29%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +2.64%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +7.62%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.994.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +11.83%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -0.01%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'92.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.37%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'16.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +36.09%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -4.18%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +10.23%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'102.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +8.52%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.704.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +8.39%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.187"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.97%  "
$ws.Range("E51").Style = "Normal"
